$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force Text format to preserve exact string representation
# (e.g. trailing zeros like "11.10", "8.00" that Excel would otherwise normalize as numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.014.24"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.579.29"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.41"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.38"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.09"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.043.49"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.925.51"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.582.33"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.93"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.166"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.00"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "461.57"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.79"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.87"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "158.74"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.08"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.70"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.29"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.635"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0538"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.05"

# Volume(1h) column (E) updates
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +7.72%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("E42").Value = "  +4.98%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -3.14%  "
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("E50").Value = "  -2.04%  "
